$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("#gratia", "Gratia"),
    @("#philib", "Philib"),
    @("#arch", "Arch"),
    @("#amou", "Amou"),
    @("#mansh", "Mansh"),
    @("#bradem", "Bradem"),
    @("#aardighe", "Aardighe"),
    @("#kar", "Kar"),
    @("#brad", "Brad"),
    @("#vvaren", "VVaren"),
    @("#aartsche-diane-ghevveldich-versteurt", "Aartsche Diane ghevveldich versteurt"),
    @("#coster", "Coster"),
    @("#edelen", "Edelen"),
    @("#aardige", "Aardige"),
    @("#lymius", "Lymius"),
    @("#amoureusje.-manshooft", "Amoureusje. Manshooft"),
    @("#2.-choor", "2. Choor"),
    @("#limius", "Limius"),
    @("#vrede", "Vrede"),
    @("#de-keyzer", "De Keyzer"),
    @("#tvvee-turcksche-iagerinnetjes", "Tvvee Turcksche Iagerinnetjes"),
    @("#aartsche-diana", "Aartsche Diana"),
    @("#doctoor", "Doctoor"),
    @("#lethea", "Lethea"),
    @("#keyz", "Keyz"),
    @("#phillib", "Phillib"),
    @("#aartsche", "Aartsche"),
    @("#keyzer", "Keyzer"),
    @("#herem", "Herem"),
    @("#1.-choor", "1. Choor"),
    @("#alderecht,-aardiges-broeder", "Alderecht, Aardiges Broeder"),
    @("#twee-iagers", "Twee Iagers"),
    @("#seege", "Seege"),
    @("#zoud", "Zoud"),
    @("#vvare", "VVare"),
    @("#sijn-broeder-vrederijck", "Sijn Broeder Vrederijck"),
    @("#chooren", "Chooren"),
    @("#amour", "Amour"),
    @("#heer", "Heer"),
    @("#moers", "Moers"),
    @("#ald", "Ald"),
    @("#mod", "Mod"),
    @("#aerts", "Aerts"),
    @("#choor", "Choor"),
    @("#zeeg", "Zeeg"),
    @("#aartsche-diane", "Aartsche Diane"),
    @("#a", "A"),
    @("#aart", "Aart"),
    @("#amara", "Amara"),
    @("#adel", "Adel"),
    @("#seeg", "Seeg"),
    @("#alder", "Alder"),
    @("#alderecht,-hereman,-manshooft,-amereusje,-en-al-de", "Alderecht, Hereman, Manshooft, Amereusje, en al de"),
    @("#romul", "Romul"),
    @("#zeege", "Zeege"),
    @("#lethe", "Lethe"),
    @("#modde", "Modde"),
    @("#amera", "Amera"),
    @("#moersg", "Moersg"),
    @("#heere", "Heere"),
    @("#aarts", "Aarts"),
    @("#soud", "Soud"),
    @("#amar", "Amar"),
    @("#moer", "Moer"),
    @("#aard", "Aard"),
    @("#wandel", "Wandel"),
    @("#aer", "Aer"),
    @("#gratian", "Gratian"),
    @("#dort", "Dort"),
    @("#phille", "Phille")
)

$url = "https://www.dbnl.org/tekst/bred001chak02_01"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 1).Value = $url
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
}
